$wb = $excel.ActiveWorkbook

# --- Sheet "findNewCar" (sheet1): update runmode (column D) values ---
$ws1 = $wb.Worksheets.Item("findNewCar")
$ws1.Range("D2").Value = "N"
$ws1.Range("D3").Value = "Y"
$ws1.Range("D5:D13").Value = "N"

# --- Sheet "carNameAndPrice" (sheet2): update runmode (column D) values ---
$ws2 = $wb.Worksheets.Item("carNameAndPrice")
$ws2.Range("D2").Value = "N"
$ws2.Range("D5:D13").Value = "N"

# --- View state: active sheet moves to findNewCar, with D13 selected ---
$ws1.Activate() | Out-Null
$ws1.Range("D13").Select() | Out-Null
